$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4 through 83 contain the per-province data; columns V and W hold the
# "Effective scale up for non-poor people" / "Effective scale up for poor
# people" values which were previously split evenly across three hazards
# (1/3 each). They are now moved to a separate demo workbook, so zero them
# out here.
$ws.Range("V4:W83").Value = 0
